# Update the "Förändrad" (changed) date column (C) for rows 2-41
# from 45204 (2023-10-05) to 45207 (2023-10-08).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2:C41").Value = 45207
